$d = $word.ActiveDocument

# Locate the paragraph that holds the "Supplementary Table 2. ..." caption text.
$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$r = $p.Range

# Make "Supplementary Table 2." bold (Word will split runs/rPr automatically).
$boldEnd = $r.Start + [int]"Supplementary Table 2.".Length
$boldRange = $d.Range($r.Start, $boldEnd)
$boldRange.Bold = 1

# Move the "_GoBack" bookmark from the now-empty preceding paragraph into the
# middle of "Neuraminidase", i.e. right after "...Results from Neura".
$bmPos = $r.Start + [int]"Supplementary Table 2. Results from Neura".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
